# Change the "Sunday" dates in column B (rows 2,4,6,...,30) to the following
# "Monday" (i.e. add one day), leaving the "Thursday" dates (odd rows)
# untouched. This mirrors the commit "changed the if from sunday to monday".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "07-02-2022"
    4  = "14-02-2022"
    6  = "21-02-2022"
    8  = "28-02-2022"
    10 = "07-03-2022"
    12 = "14-03-2022"
    14 = "21-03-2022"
    16 = "28-03-2022"
    18 = "04-04-2022"
    20 = "11-04-2022"
    22 = "18-04-2022"
    24 = "25-04-2022"
    26 = "02-05-2022"
    28 = "09-05-2022"
    30 = "16-05-2022"
}

# Rows whose new day-of-month is <= 12 are ambiguous (Excel could read them
# as mm-dd-yyyy instead of dd-mm-yyyy), so force those particular cells to
# text first to keep the literal string. Unambiguous ones (day > 12) are
# left completely alone, exactly like Excel would treat typed text.
$ambiguousRows = @(2, 10, 18, 20, 26, 28)

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("B$row")
    if ($ambiguousRows -contains $row) {
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$row]
        $cell.NumberFormat = "General"
    } else {
        $cell.Value = $updates[$row]
    }
}
